# Add maintenance cycle API: expand the AssetUsefulLife rule table with
# additional asset categories and their useful-life values (in years).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any inherited style (e.g. wrap-text) on the "LIFE" values column for the
# rows we are about to (re)populate, since the new rows use the default style.
$ws.Range("C9:C37").ClearFormats()

# Row 9 (existing row, value updated)
$ws.Range("B9").Value = "AssetCategory.TV"
$ws.Range("C9").Value = 10

# New rows 10-35: additional asset categories
$ws.Range("B10").Value = "AssetCategory.SOUND_SYSTEM"
$ws.Range("C10").Value = 12
$ws.Range("B11").Value = "AssetCategory.PROJECTOR"
$ws.Range("C11").Value = 8
$ws.Range("B12").Value = "AssetCategory.AIR_CONDITIONER"
$ws.Range("C12").Value = 15
$ws.Range("B13").Value = "AssetCategory.LIGHTING"
$ws.Range("C13").Value = 15
$ws.Range("B14").Value = "AssetCategory.AIR_PURIFIER"
$ws.Range("C14").Value = 7
$ws.Range("B15").Value = "AssetCategory.STOVE"
$ws.Range("C15").Value = 15
$ws.Range("B16").Value = "AssetCategory.MICROWAVE"
$ws.Range("C16").Value = 10
$ws.Range("B17").Value = "AssetCategory.OVEN"
$ws.Range("C17").Value = 15
$ws.Range("B18").Value = "AssetCategory.REFRIGERATOR"
$ws.Range("C18").Value = 20
$ws.Range("B19").Value = "AssetCategory.WATER_PURIFIER"
$ws.Range("C19").Value = 7
$ws.Range("B20").Value = "AssetCategory.RANGE_HOOD"
$ws.Range("C20").Value = 15
$ws.Range("B21").Value = "AssetCategory.BED"
$ws.Range("C21").Value = 15
$ws.Range("B22").Value = "AssetCategory.CAR"
$ws.Range("C22").Value = 15
$ws.Range("B23").Value = "AssetCategory.WARDROBE"
$ws.Range("C23").Value = 20
$ws.Range("B24").Value = "AssetCategory.FAN"
$ws.Range("C24").Value = 10
$ws.Range("B25").Value = "AssetCategory.LAMP"
$ws.Range("C25").Value = 15
$ws.Range("B26").Value = "AssetCategory.SHOWER"
$ws.Range("C26").Value = 15
$ws.Range("B27").Value = "AssetCategory.BATHTUB"
$ws.Range("C27").Value = 25
$ws.Range("B28").Value = "AssetCategory.SINK"
$ws.Range("C28").Value = 20
$ws.Range("B29").Value = "AssetCategory.WATER_HEATER"
$ws.Range("C29").Value = 12
$ws.Range("B30").Value = "AssetCategory.EXHAUST_FAN"
$ws.Range("C30").Value = 15
$ws.Range("B31").Value = "AssetCategory.WASHING_MACHINE"
$ws.Range("C31").Value = 12
$ws.Range("B32").Value = "AssetCategory.DRYER"
$ws.Range("C32").Value = 12
$ws.Range("B33").Value = "AssetCategory.CLOTHES_RACK"
$ws.Range("C33").Value = 10
$ws.Range("B34").Value = "AssetCategory.IRON"
$ws.Range("C34").Value = 8
$ws.Range("B35").Value = "AssetCategory.CHAIR"
$ws.Range("C35").Value = 15

# Rows 36-37: previously rows 10-11 (LAPTOP/MOBILE_PHONE), moved down and
# re-written without the inherited style.
$ws.Range("B36").Value = "AssetCategory.LAPTOP"
$ws.Range("C36").Value = 5
$ws.Range("B37").Value = "AssetCategory.MOBILE_PHONE"
$ws.Range("C37").Value = 4

# Update the view state to match what was left selected/scrolled in the edit.
$ws.Range("B15").Select()
$excel.ActiveWindow.ScrollRow = 3
